# Refresh the cryptos list: updated Price (D) / Volume(1h) (E) figures for
# the latest data pull, plus two rank swaps where the underlying coin data
# moved (Celestia <-> ARBITRUM, BEAM <-> TheGraph).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price text such as "350.43" would otherwise be auto-coerced to a number
# by Excel (losing e.g. trailing zeros like "1.00" -> 1). Prefix those with
# a quote so they stay text, same as the source data's inline strings.
function Set-PriceCell {
    param($row, $value)
    if ($null -eq $value) { return }
    if ($value -match '^[+-]?[0-9]*\.?[0-9]+$') {
        $ws.Range("D$row").Value = "'" + $value
    } else {
        $ws.Range("D$row").Value = $value
    }
}

function Set-VolumeCell {
    param($row, $value)
    if ($null -eq $value) { return }
    $ws.Range("E$row").Value = $value
}

Set-PriceCell  2 "51.644.26"
Set-VolumeCell 2 "  -1.07%  "

Set-PriceCell  3 "2.914.68"
Set-VolumeCell 3 "  -0.51%  "

Set-VolumeCell 4 "  +0.06%  "

Set-PriceCell  5 "350.43"
Set-VolumeCell 5 "  -1.09%  "

Set-PriceCell  6 "106.55"
Set-VolumeCell 6 "  -5.94%  "

Set-PriceCell  7 "0.552"
Set-VolumeCell 7 "  -1.40%  "

Set-VolumeCell 8 "  -0.01%  "

Set-PriceCell  9 "0.608"
Set-VolumeCell 9 "  -2.69%  "

Set-PriceCell  10 "37.52"
Set-VolumeCell 10 "  -5.43%  "

Set-VolumeCell 11 "  +0.92%  "

Set-PriceCell  12 "0.0851"
Set-VolumeCell 12 "  -3.62%  "

Set-PriceCell  13 "18.88"

Set-PriceCell  14 "3.389.11"
Set-VolumeCell 14 "  +0.27%  "

Set-PriceCell  15 "7.62"
Set-VolumeCell 15 "  -2.07%  "

Set-PriceCell  16 "2.922.28"
Set-VolumeCell 16 "  +0.05%  "

Set-PriceCell  17 "0.964"
Set-VolumeCell 17 "  -2.61%  "

Set-PriceCell  18 "51.596.20"
Set-VolumeCell 18 "  -1.28%  "

Set-PriceCell  19 "3.46"
Set-VolumeCell 19 "  +4.56%  "

Set-PriceCell  20 "7.32"
Set-VolumeCell 20 "  -3.75%  "

Set-PriceCell  21 "13.37"
Set-VolumeCell 21 "  -5.70%  "

Set-VolumeCell 22 "  -2.20%  "

Set-PriceCell  23 "68.80"
Set-VolumeCell 23 "  -3.46%  "

Set-PriceCell  24 "260.99"
Set-VolumeCell 24 "  -3.45%  "

Set-PriceCell  25 "2.69"
Set-VolumeCell 25 "  -4.48%  "

Set-PriceCell  26 "7.64"
Set-VolumeCell 26 "  +8.07%  "

Set-PriceCell  27 "0.171"
Set-VolumeCell 27 "  -4.99%  "

Set-PriceCell  28 "26.42"
Set-VolumeCell 28 "  -1.97%  "

Set-PriceCell  29 "1.00"

Set-VolumeCell 30 "  -0.58%  "

Set-PriceCell  31 "10.19"
Set-VolumeCell 31 "  -4.56%  "

Set-VolumeCell 32 "  +1.43%  "

Set-VolumeCell 33 "  -4.93%  "

Set-PriceCell  34 "35.45"
Set-VolumeCell 34 "  -4.78%  "

Set-VolumeCell 35 "  -4.43%  "

Set-VolumeCell 37 "  -6.50%  "

Set-PriceCell  38 "3.10"
Set-VolumeCell 38 "  -8.34%  "

# Rows 39/40: ARBITRUM now outranks Celestia, swap their data.
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-PriceCell  39 "1.94"
Set-VolumeCell 39 "  -5.53%  "

$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-PriceCell  40 "17.54"
Set-VolumeCell 40 "  -6.70%  "

Set-PriceCell  41 "2.64"
Set-VolumeCell 41 "  -2.87%  "

Set-VolumeCell 42 "  -2.08%  "

Set-PriceCell  43 "22.33"
Set-VolumeCell 43 "  -3.27%  "

Set-PriceCell  44 "119.44"
Set-VolumeCell 44 "  +2.78%  "

Set-PriceCell  45 "2.14"
Set-VolumeCell 45 "  -2.24%  "

Set-PriceCell  46 "2.090.25"
Set-VolumeCell 46 "  -4.63%  "

Set-PriceCell  47 "3.30"
Set-VolumeCell 47 "  -6.99%  "

Set-PriceCell  48 "2.30"
Set-VolumeCell 48 "  -8.83%  "

# Rows 49/50: TheGraph now outranks BEAM, swap their data.
$ws.Range("B49").Value = "TheGraph"
$ws.Range("C49").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-PriceCell  49 "0.237"
Set-VolumeCell 49 "  -5.86%  "

$ws.Range("B50").Value = "BEAM"
$ws.Range("C50").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
Set-PriceCell  50 "0.0339"
Set-VolumeCell 50 "  -3.13%  "

Set-VolumeCell 51 "  -7.02%  "
